# "added mid term course on web"
# Sheet1 tracks three courses (TPL = column C, SRE = column D, ALGO = column E),
# one row per "week". This edit:
#  - adds "Slides:" lines to the three existing ALGO week entries (E7, E8->moved, E9->moved)
#  - inserts a new "Mid Term Exam" row entry for each of the three courses
#    (TPL -> C8, SRE -> D9, ALGO -> E10)
#  - grows rows 9 & 10 to fit the extra wrapped text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- E7: ALGO Week 5 gains a Slides line -------------------------------
$e7 = @"
Topic: name- Week 5,lectures- 2 Lectures, duration- 01:16;
Video: link- https://drive.google.com/file/d/1t6VMFiI4610ULxoovPeFNuAnRKYwKosu/preview, name- AD&AA Week#5 Part 1.mp4, duration- 00:25;
Video: link- https://drive.google.com/file/d/1s_RbEVFMh9pGqPTkgfURH9Kx8nGcK3GI/preview, name- AD&AA Week#5 Part 2.mp4, duration- 00:51;
Slides: slide- lec3a.ppt
"@
$ws.Cells.Item(7, 5).Value = $e7

# --- C8: TPL Mid Term Exam (new content; cell was previously blank) ----
$c8 = @"
Topic: name- Mid Term Exam, lectures- , duration- 11 Nov | 15:00-16:30;
Slides: slide- Chapter 1.ppt;
Slides: slide- Chapter 2.ppt;
Slides: slide- Chapter 3.ppt;
Slides: slide- Chapter 5.pptx;
Books: name- Concepts of Programming Languages 11th Ed, link- Concepts of Programming Languages 11th Ed.pdf;
"@
$ws.Cells.Item(8, 3).Value = $c8
$ws.Cells.Item(8, 3).Font.Bold = $true
$ws.Cells.Item(8, 3).WrapText = $true
$ws.Cells.Item(8, 3).HorizontalAlignment = -4131
$ws.Cells.Item(8, 3).VerticalAlignment = -4108

# --- E8: ALGO Week 6 gains a Slides line --------------------------------
$e8 = @"
Topic: name- Week 6,lectures- 2 Lectures, duration- 01:42;
Video: link- https://drive.google.com/file/d/1rxg0O9yrSbrJpaYiql0F-P8u4Y2ZtDuc/preview, name- AD&AA Week#6 Part 1.mp4, duration- 00:38;
Video: link- https://drive.google.com/file/d/13r2bglsKWo9HeYgPDRZidU9CV5oJW6xu/preview, name- AD&AA Week#6 Part 2.mp4, duration- 01:04;
Slides: slide- lec3b.pptx;
"@
$ws.Cells.Item(8, 5).Value = $e8

# --- D9: SRE Mid Term Exam (new content; cell was previously blank) ----
$d9 = @"
Topic: name- Mid Term Exam, lectures- , duration- 12 Nov | 15:00-16:30;
Slides: slide- Lecture Slide_1.pptx;
Slides: slide- Lecture Slide_2.pptx;
Slides: slide- Lecture Slide_3.pptx;
Slides: slide- Lecture Slide_4.pptx;
Books: name- Requirements Engineering Fundamentals, link- Requirements Engineering Fundamentals A Study Guide for the Certified Professional for Requirements Engineering Exam.pdf;
"@
$ws.Cells.Item(9, 4).Value = $d9

# --- E9: ALGO Week 7 gains a Slides line --------------------------------
$e9 = @"
Topic: name- Week 7,lectures- 2 Lectures, duration- 01:41;
Video: link- https://drive.google.com/file/d/1iN0XDjE-kn0od8h8e3azwQFUaFEwDKEt/preview, name- AD&AA Week#7 Part 1.mp4, duration- 00:57;
Video: link- https://drive.google.com/file/d/130Cyy5mWSzIPNuH8qWlAk8wQcO5H_7o4/preview, name- AD&AA Week#7 Part 2.mp4, duration- 00:44;
Slides: slide- lec4.pptx;
"@
$ws.Cells.Item(9, 5).Value = $e9

# --- E10: ALGO Mid Term Exam (new content; cell was previously blank) ---
$e10 = @"
Topic: name- Mid Term Exam, lectures- , duration- 8 Nov | 16:00-17:30;
Assignment: name- Mid Term Exam Pattern, img- /ALGO/paperpattern.png;
Slides: slide- lec1.pptx;
Slides: slide- lec2.pdf;
Slides: slide- lec3b.pptx;
"@
$ws.Cells.Item(10, 5).Value = $e10
$ws.Cells.Item(10, 5).Font.Bold = $true
$ws.Cells.Item(10, 5).WrapText = $true
$ws.Cells.Item(10, 5).HorizontalAlignment = -4131
$ws.Cells.Item(10, 5).VerticalAlignment = -4108

# --- Row heights grow to fit the newly-added lines ----------------------
$ws.Rows.Item(9).RowHeight = 170
$ws.Rows.Item(10).RowHeight = 119

# --- Selection moves to the newly added Mid Term cell for ALGO ----------
$ws.Activate()
$ws.Range("E10").Select()
